$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.770.95"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "3.317.96"
$ws.Range("E3").Value = "  +5.17%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'603.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.84%  "
$ws.Range("E6").Value = "  +3.03%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.315.42"
$ws.Range("E8").Value = "  +5.08%  "
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("E10").Value = "  +2.92%  "
$ws.Range("D11").Value = "'5.50"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.04%  "
$ws.Range("E12").Value = "  +2.35%  "
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "'34.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("D15").Value = "3.864.56"
$ws.Range("E15").Value = "  +5.23%  "
$ws.Range("D16").Value = "'0.120"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "3.320.54"
$ws.Range("E17").Value = "  +5.28%  "
$ws.Range("D18").Value = "63.847.81"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("E19").Value = "  +3.38%  "
$ws.Range("D20").Value = "'480.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("E22").Value = "  +5.16%  "
$ws.Range("D23").Value = "'8.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.60%  "
$ws.Range("D24").Value = "'13.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.63%  "
$ws.Range("D25").Value = "'84.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +2.66%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "'7.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.66%  "
$ws.Range("D30").Value = "'8.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.36%  "
$ws.Range("E31").Value = "  +3.73%  "
$ws.Range("D32").Value = "'28.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.95%  "
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("D35").Value = "'1.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.41%  "
$ws.Range("E36").Value = "  +4.99%  "
$ws.Range("D37").Value = "'52.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("D38").Value = "0.0₃0740"
$ws.Range("E38").Value = "  +5.58%  "
$ws.Range("E39").Value = "  +3.66%  "
$ws.Range("D40").Value = "'433.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.98%  "
$ws.Range("D41").Value = "3.104.43"
$ws.Range("E41").Value = "  +4.97%  "
$ws.Range("D42").Value = "'0.121"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.20%  "
$ws.Range("D43").Value = "'2.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("E45").Value = "  +2.43%  "
$ws.Range("E46").Value = "  +4.93%  "
$ws.Range("D47").Value = "'37.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +16.48%  "
$ws.Range("D48").Value = "'26.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.70%  "
$ws.Range("D50").Value = "'2.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.78%  "
$ws.Range("E51").Value = "  +0.81%  "
